# Added Indian MF 1st Stab
# Insert 9 new date columns (Jun_16 .. Sep_08) right after column A,
# pushing the existing date columns (previously B:S) to the right (K:AB).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 9 new columns before column B (old B becomes K, old S becomes AB)
$ws.Range("B1:J1").EntireColumn.Insert()

# 2. Populate the new header cells (row 1) with the new date labels,
#    most-recent-first just like the existing header row.
$ws.Range("B1").Value2 = "Sep_08"
$ws.Range("C1").Value2 = "Aug_25"
$ws.Range("D1").Value2 = "Aug_04"
$ws.Range("E1").Value2 = "Jul_23"
$ws.Range("F1").Value2 = "Jul_17"
$ws.Range("G1").Value2 = "Jul_07"
$ws.Range("H1").Value2 = "Jun_30"
$ws.Range("I1").Value2 = "Jun_24"
$ws.Range("J1").Value2 = "Jun_16"

# 3. Populate the new data cells with the "UN" placeholder rating used
#    throughout the rest of the sheet, matching each row's existing extent.
$ws.Range("B2:J29").Value2 = "UN"
$ws.Range("B30:J31").Value2 = "UN"
$ws.Range("B32:J33").Value2 = "UN"

# 4. Extend the existing column width formatting (8.0 chars, same as the
#    pre-existing C:S columns) across the newly shifted column range C:AB.
$ws.Range("C1:AB1").EntireColumn.ColumnWidth = 7.14
